$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before J. This shifts the existing J/K/L columns
#    (Photo taken?, External Profile Updated?, Tested?) one place to the
#    right, turning them into K/L/M, and creates a blank new column J.
[void]$ws.Columns("J:J").Insert()

# 2. Give the new column its header text.
$ws.Range("J1").Value = "IsExternal?"

# 3. Populate the new "IsExternal?" boolean column for every data row.
#    Default everybody to FALSE (SSW staff)...
$ws.Range("J2:J30").Value = $false
# ...except the two non-SSW guest profiles (Andrew Coates & Samantha Coates)
# which are external and should be TRUE.
$ws.Range("J29").Value = $true
$ws.Range("J30").Value = $true

# 4. Match the cell formatting/border style used by the rest of the table
#    (style of the data cells, e.g. K2) for the whole new column.
$ws.Range("K2").Copy()
$ws.Range("J2:J30").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Give the new column a sensible width (matching its neighbour) instead of
# the sheet default, similar to the other wide data columns H/I.
$ws.Columns("J").ColumnWidth = $ws.Columns("I").ColumnWidth

# 5. Re-point the worksheet AutoFilter so it covers the new column M
#    (A1:L28 -> A1:M28). Turn it off first, temporarily insert a blank
#    separator row below the filtered table so the engine doesn't expand
#    the detected range into the extra profile rows (29/30), reapply the
#    filter on the correct range, then remove the helper row again.
$ws.AutoFilterMode = $false
[void]$ws.Rows("29:29").Insert()
[void]$ws.Range("A1:M28").AutoFilter()
[void]$ws.Rows("29:29").Delete()

# 6. Update the workbook-level _FilterDatabase defined name to match.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$M`$28"
  }
}

# 7. Leave the selection on the last cell touched, like a real editing
#    session would.
[void]$ws.Range("I30").Select()
